$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1, J1 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (style) from H1 (an existing header cell) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data values for columns I and J, rows 2-32 ---
$iValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 6; 27 = 1; 28 = 1;
    29 = 1; 30 = 4; 31 = 1; 32 = 3
}
$jValues = @{
    2 = 4; 3 = 3; 4 = 5; 5 = 6; 6 = 6; 7 = 6; 8 = 8; 9 = 8; 10 = 6;
    11 = 7; 12 = 4; 13 = 6; 14 = 5; 15 = 7; 16 = 1; 17 = 6; 18 = 5; 19 = 4;
    20 = 5; 21 = 6; 22 = 2; 23 = 4; 24 = 5; 25 = 4; 26 = 9; 27 = 5; 28 = 5;
    29 = 3; 30 = 6; 31 = 2; 32 = 4
}

for ($row = 2; $row -le 32; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
